$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 136, shifting existing rows 136:149 down to 137:150
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with the new record
$ws.Cells.Item(136, 1).Value = 4
$ws.Cells.Item(136, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(136, 3).Value = "Los Lagos"
$ws.Cells.Item(136, 4).Value = 44543
$ws.Cells.Item(136, 5).Value = 10
$ws.Cells.Item(136, 6).Value = 100112028
$ws.Cells.Item(136, 7).Value = "Sandia"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Tercera"
$ws.Cells.Item(136, 10).Value = 400
$ws.Cells.Item(136, 11).Value = 2500
$ws.Cells.Item(136, 12).Value = 2500
$ws.Cells.Item(136, 13).Value = 2500
$ws.Cells.Item(136, 14).Value = "`$/unidad"
$ws.Cells.Item(136, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(136, 16).Value = 2500
$ws.Cells.Item(136, 17).Value = 1
$ws.Cells.Item(136, 18).Value = "Hortaliza"
